$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "time_taken" header in column F, matching the style of the
# existing header row (B1:E1) by copying formats from E1.
$ws.Range("F1").Value = "time_taken"
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Populate time_taken values for each data row (2-64).
$ws.Range("F2").Value = "2021-10-05 13:38:39.202385"
$ws.Range("F3").Value = "2021-10-05 13:38:39.202396"
$ws.Range("F4").Value = "2021-10-05 13:38:39.202400"
$ws.Range("F5").Value = "2021-10-05 13:38:39.202402"
$ws.Range("F6").Value = "2021-10-05 13:38:39.202405"
$ws.Range("F7").Value = "2021-10-05 13:38:39.202408"
$ws.Range("F8").Value = "2021-10-05 13:38:39.202410"
$ws.Range("F9").Value = "2021-10-05 13:38:39.202413"
$ws.Range("F10").Value = "2021-10-05 13:38:39.202415"
$ws.Range("F11").Value = "2021-10-05 13:38:39.202418"
$ws.Range("F12").Value = "2021-10-05 13:38:39.202420"
$ws.Range("F13").Value = "2021-10-05 13:38:39.202423"
$ws.Range("F14").Value = "2021-10-05 13:38:39.202425"
$ws.Range("F15").Value = "2021-10-05 13:38:39.202428"
$ws.Range("F16").Value = "2021-10-05 13:38:39.202430"
$ws.Range("F17").Value = "2021-10-05 13:38:39.202432"
$ws.Range("F18").Value = "2021-10-05 13:38:39.202435"
$ws.Range("F19").Value = "2021-10-05 13:38:39.202438"
$ws.Range("F20").Value = "2021-10-05 13:38:39.202440"
$ws.Range("F21").Value = "2021-10-05 13:38:39.202442"
$ws.Range("F22").Value = "2021-10-05 13:38:39.202445"
$ws.Range("F23").Value = "2021-10-05 13:38:39.202447"
$ws.Range("F24").Value = "2021-10-05 13:38:39.202450"
$ws.Range("F25").Value = "2021-10-05 13:38:39.202452"
$ws.Range("F26").Value = "2021-10-05 13:38:39.202455"
$ws.Range("F27").Value = "2021-10-05 13:38:39.202458"
$ws.Range("F28").Value = "2021-10-05 13:38:39.202461"
$ws.Range("F29").Value = "2021-10-05 13:38:39.202463"
$ws.Range("F30").Value = "2021-10-05 13:38:39.202466"
$ws.Range("F31").Value = "2021-10-05 13:38:39.202468"
$ws.Range("F32").Value = "2021-10-05 13:38:39.202470"
$ws.Range("F33").Value = "2021-10-05 13:38:39.202473"
$ws.Range("F34").Value = "2021-10-05 13:38:39.202476"
$ws.Range("F35").Value = "2021-10-05 13:38:39.202478"
$ws.Range("F36").Value = "2021-10-05 13:38:39.202481"
$ws.Range("F37").Value = "2021-10-05 13:38:39.202483"
$ws.Range("F38").Value = "2021-10-05 13:38:39.202485"
$ws.Range("F39").Value = "2021-10-05 13:38:39.202488"
$ws.Range("F40").Value = "2021-10-05 13:38:39.202490"
$ws.Range("F41").Value = "2021-10-05 13:38:39.202493"
$ws.Range("F42").Value = "2021-10-05 13:38:39.202496"
$ws.Range("F43").Value = "2021-10-05 13:38:39.202498"
$ws.Range("F44").Value = "2021-10-05 13:38:39.202501"
$ws.Range("F45").Value = "2021-10-05 13:38:39.202503"
$ws.Range("F46").Value = "2021-10-05 13:38:39.202505"
$ws.Range("F47").Value = "2021-10-05 13:38:39.202508"
$ws.Range("F48").Value = "2021-10-05 13:38:39.202510"
$ws.Range("F49").Value = "2021-10-05 13:38:39.202513"
$ws.Range("F50").Value = "2021-10-05 13:38:39.202515"
$ws.Range("F51").Value = "2021-10-05 13:38:39.202518"
$ws.Range("F52").Value = "2021-10-05 13:38:39.202520"
$ws.Range("F53").Value = "2021-10-05 13:38:39.202522"
$ws.Range("F54").Value = "2021-10-05 13:38:39.202525"
$ws.Range("F55").Value = "2021-10-05 13:38:39.202528"
$ws.Range("F56").Value = "2021-10-05 13:38:39.202531"
$ws.Range("F57").Value = "2021-10-05 13:38:39.202533"
$ws.Range("F58").Value = "2021-10-05 13:38:39.202535"
$ws.Range("F59").Value = "2021-10-05 13:38:39.202538"
$ws.Range("F60").Value = "2021-10-05 13:38:39.202540"
$ws.Range("F61").Value = "2021-10-05 13:38:39.202543"
$ws.Range("F62").Value = "2021-10-05 13:38:39.202545"
$ws.Range("F63").Value = "2021-10-05 13:38:39.202548"
$ws.Range("F64").Value = "2021-10-05 13:38:39.202550"

